$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: DKS / Dekstop / Desktop Computer / eng
$ws.Range("A2").Value = "DKS"
$ws.Range("B2").Value = "Dekstop"
$ws.Range("C2").Value = "Desktop Computer"
$ws.Range("D2").Value = "eng"

# Row 3: DKS / الحاسوب / أجهزة الكمبيوتر المكتبية / ara
$ws.Range("A3").Value = "DKS"
$ws.Range("B3").Value = "الحاسوب"
$ws.Range("C3").Value = "أجهزة الكمبيوتر المكتبية"
$ws.Range("D3").Value = "ara"

# Row 4: DKS / Ordinateur / Ordinateurs de bureau / fra
$ws.Range("A4").Value = "DKS"
$ws.Range("B4").Value = "Ordinateur"
$ws.Range("C4").Value = "Ordinateurs de bureau"
$ws.Range("D4").Value = "fra"

[void]$ws.Range("D10").Select()

$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
